$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 values: Date, 24 hourly prices (B:Y), Price_Daily_Avg (Z), Slot_4h_max (AA),
# Slot_4h_price (AB), Slot_2h_frist (AC), Slot_2h_frist_price (AD), Slot_2h_second (AE),
# Slot_2h_second_price (AF), Slot_min_price (AG)

$ws.Range("A2").Value = 45873

$ws.Range("B2").Value = 80.98999999999999
$ws.Range("C2").Value = 80.40000000000001
$ws.Range("D2").Value = 80
$ws.Range("E2").Value = 80.40000000000001
$ws.Range("F2").Value = 85.27
$ws.Range("G2").Value = 85.27
$ws.Range("H2").Value = 95.67
$ws.Range("I2").Value = 106.23
$ws.Range("J2").Value = 97.81
$ws.Range("K2").Value = 52.57
$ws.Range("L2").Value = 30
$ws.Range("M2").Value = 13.75
$ws.Range("N2").Value = 6.76
$ws.Range("O2").Value = 4.31
$ws.Range("P2").Value = 4.31
$ws.Range("Q2").Value = 4.31
$ws.Range("R2").Value = 5.79
$ws.Range("S2").Value = 26.35
$ws.Range("T2").Value = 38.94
$ws.Range("U2").Value = 70.77
$ws.Range("V2").Value = 103.25
$ws.Range("W2").Value = 168.02
$ws.Range("X2").Value = 138.8
$ws.Range("Y2").Value = 115.04
$ws.Range("Z2").Value = 65.63

$ws.Range("AA2").Value = "20h-24h"
$ws.Range("AB2").Value = 131.28
$ws.Range("AC2").Value = "20h-22h"
$ws.Range("AD2").Value = 135.64
$ws.Range("AE2").Value = "22h-24h"
$ws.Range("AF2").Value = 126.92
$ws.Range("AG2").Value = "9h-18h"
